$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Status column (F) values for the rows whose status changed
$ws.Range("F3").Value = "Completed"
$ws.Range("F4").Value = "In progress"
$ws.Range("F9").Value = "In progress"
$ws.Range("F14").Value = "In progress"

# Scroll the sheet so row 3 becomes the top visible row, then move the
# active selection to G14
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("G14").Select()
